$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 1499
$ws.Range("I31").Value = 1499
$ws.Range("K31").Value = 4497
$ws.Range("M31").Value = -4267
$ws.Range("H74").Value = 3440
$ws.Range("I74").Value = 3484.4443
$ws.Range("K74").Value = 3484.4443
$ws.Range("M74").Value = -2548.4443
$ws.Range("H77").Value = 3440
$ws.Range("I77").Value = 3484.4443
$ws.Range("K77").Value = 17422.2215
$ws.Range("M77").Value = -12742.2215
$ws.Range("H86").Value = 5264835.5
$ws.Range("I86").Value = 7693459.5
$ws.Range("J86").Value = 2817
$ws.Range("K86").Value = 7693459.5
$ws.Range("L86").Value = 2817
$ws.Range("M86").Value = -7692336.5
$ws.Range("N86").Value = -5063
$ws.Range("H89").Value = 5264835.5
$ws.Range("I89").Value = 7693459.5
$ws.Range("J89").Value = 2817
$ws.Range("K89").Value = 38467297.5
$ws.Range("L89").Value = 14085
$ws.Range("M89").Value = -38461681.5
$ws.Range("N89").Value = -25317
$ws.Range("H132").Value = 21898.844
$ws.Range("I132").Value = 3008.9143
$ws.Range("K132").Value = 9026.742899999999
$ws.Range("M132").Value = -6496.742899999999
$ws.Range("H137").Value = 2201935.2
$ws.Range("I137").Value = 5495353
$ws.Range("J137").Value = 6323.381
$ws.Range("K137").Value = 16486059
$ws.Range("L137").Value = 18970.143
$ws.Range("M137").Value = -16483509
$ws.Range("N137").Value = -24070.143

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3466.1538
$ws.Range("I61").Value = 2598.75
$ws.Range("J61").Value = 3851.6667
$ws.Range("K61").Value = 2598.75
$ws.Range("L61").Value = 3851.6667
$ws.Range("M61").Value = -2386.75
$ws.Range("N61").Value = -4275.6667
$ws.Range("H74").Value = 1997.825
$ws.Range("I74").Value = 1726.6897
$ws.Range("K74").Value = 1726.6897
$ws.Range("M74").Value = -852.6896999999999
$ws.Range("H77").Value = 1997.825
$ws.Range("I77").Value = 1726.6897
$ws.Range("K77").Value = 8633.448499999999
$ws.Range("M77").Value = -4265.448499999999
$ws.Range("H122").Value = 2237.8333
$ws.Range("I122").Value = 2006.75
$ws.Range("K122").Value = 6020.25
$ws.Range("M122").Value = -3570.25
$ws.Range("H136").Value = 3466.1538
$ws.Range("I136").Value = 2598.75
$ws.Range("J136").Value = 3851.6667
$ws.Range("K136").Value = 7796.25
$ws.Range("L136").Value = 11555.0001
$ws.Range("M136").Value = -5246.25
$ws.Range("N136").Value = -16655.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 3000
$ws.Range("I29").Value = 3000
$ws.Range("K29").Value = 3000
$ws.Range("M29").Value = -2711
$ws.Range("H76").Value = 23500
$ws.Range("J76").Value = 23500
$ws.Range("L76").Value = 23500
$ws.Range("N76").Value = -24130
$ws.Range("H79").Value = 23500
$ws.Range("J79").Value = 23500
$ws.Range("L79").Value = 23500
$ws.Range("N79").Value = -25684
$ws.Range("H87").Value = 25000
$ws.Range("J87").Value = 25000
$ws.Range("L87").Value = 25000
$ws.Range("N87").Value = -27496
$ws.Range("H90").Value = 25000
$ws.Range("J90").Value = 25000
$ws.Range("L90").Value = 75000
$ws.Range("N90").Value = -87480

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7757720.5
$ws.Range("I31").Value = 2096.625
$ws.Range("J31").Value = 12353646
$ws.Range("K31").Value = 2096.625
$ws.Range("L31").Value = 12353646
$ws.Range("M31").Value = -1801.625
$ws.Range("N31").Value = -12354236
$ws.Range("H34").Value = 7757720.5
$ws.Range("I34").Value = 2096.625
$ws.Range("J34").Value = 12353646
$ws.Range("K34").Value = 2096.625
$ws.Range("L34").Value = 12353646
$ws.Range("M34").Value = -1894.625
$ws.Range("N34").Value = -12354050
$ws.Range("H52").Value = 38333.332
$ws.Range("J52").Value = 38333.332
$ws.Range("L52").Value = 38333.332
$ws.Range("N52").Value = -38921.332
$ws.Range("H99").Value = 3823.5833
$ws.Range("I99").Value = 3564.6
$ws.Range("J99").Value = 4008.5715
$ws.Range("K99").Value = 3564.6
$ws.Range("L99").Value = 4008.5715
$ws.Range("M99").Value = -2066.6
$ws.Range("N99").Value = -7004.5715
$ws.Range("H126").Value = 3823.5833
$ws.Range("I126").Value = 3564.6
$ws.Range("J126").Value = 4008.5715
$ws.Range("K126").Value = 10693.8
$ws.Range("L126").Value = 12025.7145
$ws.Range("M126").Value = -8223.799999999999
$ws.Range("N126").Value = -16965.7145

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 3287572.8
$ws.Range("I11").Value = 4601201
$ws.Range("J11").Value = 3502
$ws.Range("K11").Value = 4601201
$ws.Range("L11").Value = 3502
$ws.Range("M11").Value = -4601062
$ws.Range("N11").Value = -3780
$ws.Range("H122").Value = 1660
$ws.Range("I122").Value = 1733.3334
$ws.Range("J122").Value = 1550
$ws.Range("K122").Value = 5200.0002
$ws.Range("L122").Value = 4650
$ws.Range("M122").Value = -2750.0002
$ws.Range("N122").Value = -9550
$ws.Range("H126").Value = 50007220
$ws.Range("I126").Value = 100011740
$ws.Range("J126").Value = 2699.8
$ws.Range("K126").Value = 300035220
$ws.Range("L126").Value = 8099.400000000001
$ws.Range("M126").Value = -300032750
$ws.Range("N126").Value = -13039.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2591.2727
$ws.Range("I122").Value = 2550.4
$ws.Range("K122").Value = 7651.200000000001
$ws.Range("M122").Value = -5201.200000000001
$ws.Range("H132").Value = 2783.9395
$ws.Range("I132").Value = 1958.9
$ws.Range("J132").Value = 4053.2307
$ws.Range("K132").Value = 5876.700000000001
$ws.Range("L132").Value = 12159.6921
$ws.Range("M132").Value = -3346.700000000001
$ws.Range("N132").Value = -17219.6921
$ws.Range("H136").Value = 2566.0557
$ws.Range("I136").Value = 1891.5385
$ws.Range("J136").Value = 4319.8
$ws.Range("K136").Value = 5674.6155
$ws.Range("L136").Value = 12959.4
$ws.Range("M136").Value = -3124.6155
$ws.Range("N136").Value = -18059.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1300207
$ws.Range("I122").Value = 1505297.5
$ws.Range("J122").Value = 1300
$ws.Range("K122").Value = 4515892.5
$ws.Range("L122").Value = 3900
$ws.Range("M122").Value = -4513442.5
$ws.Range("N122").Value = -8800
